$d = $word.ActiveDocument
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParaXml($para, $innerXml) {
    $xml = "<w:p $ns>$innerXml</w:p>"
    $para.Range.InsertXML($xml)
}

function Set-CellXml($cell, $innerXml) {
    $xml = "<w:p $ns>$innerXml</w:p>"
    $para = $cell.Range.Paragraphs.Item(1)
    $para.Range.InsertXML($xml)
}

# ---------------------------------------------------------------------
# 1) Remove the blank paragraph between the title "Työaikaraportti" and
#    the "Ryhmä: 2" paragraph.
# ---------------------------------------------------------------------
$d.Paragraphs.Item(2).Range.Delete()

# ---------------------------------------------------------------------
# 2) "Ryhmä: 2" paragraph -> move the space: " 2" becomes "2 "
#    (keep the two runs, just swap which side the space sits on)
# ---------------------------------------------------------------------
$ryhmaPara = $d.Paragraphs.Item(2)
Set-ParaXml $ryhmaPara '<w:pPr><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Ryhmä:</w:t></w:r><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">2 </w:t></w:r>'

# ---------------------------------------------------------------------
# 3) "Henkilö: Joona Jalonen" -> "Henkilö: Iiro Anttila"
# ---------------------------------------------------------------------
$henkiloPara = $d.Paragraphs.Item(3)
Set-ParaXml $henkiloPara '<w:pPr><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">Henkilö: </w:t></w:r><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Iiro Anttila</w:t></w:r>'

# ---------------------------------------------------------------------
# 4) Wrap "Pvm" table header text with gramStart/gramEnd proofErr marks
# ---------------------------------------------------------------------
$table = $d.Tables.Item(1)
$pvmCell = $table.Rows.Item(1).Cells.Item(1)
Set-CellXml $pvmCell '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Pvm</w:t></w:r><w:proofErr w:type="gramEnd"/>'

# ---------------------------------------------------------------------
# 5) Fill in the work-log rows. There are currently 5 blank data rows
#    (rows 3-7) right after the header row (row 1) and the first blank
#    row that stays blank (row 2). The edit fills those 5 rows with
#    content and adds one extra row, for 6 filled rows total.
# ---------------------------------------------------------------------
$beforeRow = $table.Rows.Item(8)
$table.Rows.Add($beforeRow) | Out-Null

$rowsData = @(
    @{
        col1 = '<w:r><w:t>8</w:t></w:r><w:r><w:t>.3</w:t></w:r>'
        col2 = '<w:r><w:t>3</w:t></w:r>'
        col3 = '<w:r><w:t>Koodin suunnittelu ja ideoiden kertaaminen ryhmäläisen kanssa</w:t></w:r>'
    },
    @{
        col1 = '<w:r><w:t>8.3</w:t></w:r>'
        col2 = '<w:r><w:t>4</w:t></w:r>'
        col3 = '<w:r><w:t>Koodin rakennuksen aloittaminen</w:t></w:r>'
    },
    @{
        col1 = '<w:r><w:t>9.3</w:t></w:r>'
        col2 = '<w:r><w:t>5</w:t></w:r>'
        col3 = '<w:r><w:t>koodin rakentaminen</w:t></w:r>'
    },
    @{
        col1 = '<w:r><w:t>10.3</w:t></w:r>'
        col2 = '<w:r><w:t>4</w:t></w:r>'
        col3 = '<w:r><w:t>koodin viimeisteleminen</w:t></w:r>'
    },
    @{
        col1 = '<w:r><w:t>17.3</w:t></w:r>'
        col2 = '<w:r><w:t>2</w:t></w:r>'
        col3 = '<w:r><w:t xml:space="preserve">Opettelin käyttämään </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>classe</w:t></w:r><w:r><w:t>ja</w:t></w:r><w:proofErr w:type="spellEnd"/>'
    },
    @{
        col1 = '<w:r><w:t>18.3</w:t></w:r>'
        col2 = '<w:r><w:t>3</w:t></w:r>'
        col3 = '<w:r><w:t>Korjasin koodia</w:t></w:r><w:r><w:t xml:space="preserve"> (en saanut valmiiksi, joten toinen ryhmäläinen hoiti sen loppuun)</w:t></w:r>'
    }
)

for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $row = $table.Rows.Item(3 + $i)
    $data = $rowsData[$i]
    Set-CellXml $row.Cells.Item(1) $data.col1
    Set-CellXml $row.Cells.Item(2) $data.col2
    Set-CellXml $row.Cells.Item(3) $data.col3
}

# ---------------------------------------------------------------------
# 6) Fill in "21" (bold) in the bottom-right "Yhteensä" total cell
# ---------------------------------------------------------------------
$lastRow = $table.Rows.Item($table.Rows.Count)
$totalCell = $lastRow.Cells.Item($lastRow.Cells.Count)
Set-CellXml $totalCell '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>21</w:t></w:r>'

Write-Host "All edits applied"
